$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Note: cell A16's label "SGD" is a brand-new shared string (distinct from the
# pre-existing "SGD " with a trailing space used in row 3), so it is written
# first to mirror the order new strings were introduced by the original author.
$ws.Range("A16").Value = "SGD"

# --- Complete row 9 (A9/B9 already contain "KNN" / "No ") ---
$ws.Range("C9").Value = "83.264s"
$ws.Range("D9").Value = 159
$ws.Range("E9").Value = 0.96389999999999998
$ws.Range("F9").Value = 2
$ws.Range("G9").Value = 0.0625

# --- New bold section header in row 11 ---
$ws.Range("A11").Value = "More Filtration"
$ws.Range("A11").Font.Bold = $true

# --- New "More Filtration" results table, rows 12-16 ---
# Row 12: LogReg
$ws.Range("A12").Value = "LogReg"
$ws.Range("B12").Value = "Yes"
$ws.Range("C12").Value = "1.139s"
$ws.Range("D12").Value = 146
$ws.Range("E12").Value = 0.88480000000000003
$ws.Range("F12").Value = 19
$ws.Range("G12").Value = 0.59375

# Row 13: KNN
$ws.Range("A13").Value = "KNN"
$ws.Range("B13").Value = "Yes"
$ws.Range("C13").Value = "1.074s"
$ws.Range("D13").Value = 127
$ws.Range("E13").Value = 0.76968999999999999
$ws.Range("F13").Value = 21
$ws.Range("G13").Value = 0.65625

# Row 14: SVM
$ws.Range("A14").Value = "SVM "
$ws.Range("B14").Value = "Yes "
$ws.Range("C14").Value = "1.119s"
$ws.Range("D14").Value = 126
$ws.Range("E14").Value = 0.76359999999999995
$ws.Range("F14").Value = 26
$ws.Range("G14").Value = 0.8125

# Row 15: LSTM
$ws.Range("A15").Value = "LSTM"
$ws.Range("B15").Value = "Yes"
$ws.Range("C15").Value = "93.51s"
$ws.Range("D15").Value = 165
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0

# Row 16: SGD (remaining cells; A16 was already set above)
$ws.Range("B16").Value = "Yes"
$ws.Range("C16").Value = "24.358s"
$ws.Range("D16").Value = 144
$ws.Range("E16").Value = 0.87270000000000003
$ws.Range("F16").Value = 22
$ws.Range("G16").Value = 0.6875

# --- Update selection to mirror the final authoring position ---
$ws.Range("G16").Select() | Out-Null
